# Add a new worksheet "ShopForDSLRs" at the end of the workbook, matching
# the layout/content of the other "search result" sheets (Search, ViewProduct):
#   A2 -> shared text "0"
#   B2 -> shared text "1"
#   C2 -> the new search-term string, formatted like the other sheets' last
#         (Menlo-font) column.

$wb = $excel.ActiveWorkbook

# Template cell whose number format / font we want to reuse for column C
# (Menlo font, general number format -> matches style index already present
# in the workbook, so no new style gets created).
$templateCell = $wb.Worksheets.Item("ViewProduct").Range("F2")

# Insert the new sheet after the last existing sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "ShopForDSLRs"

# A2 / B2 hold the same "0" / "1" flag strings used on the ViewProduct sheet.
# Force text number format first so Excel stores them as shared strings
# instead of numbers (matching the source workbook's cell typing).
$newSheet.Range("A2").NumberFormat = "@"
$newSheet.Range("A2").Value = "0"
$newSheet.Range("B2").NumberFormat = "@"
$newSheet.Range("B2").Value = "1"

# C2 holds the new search query text.
$newSheet.Range("C2").Value = "Canon Digital SLR Cameras between `$350.00 and `$750.00"

# Match C2's formatting (Menlo font) to the analogous column on the other
# sheets by copying formats from the template cell.
$templateCell.Copy() | Out-Null
$newSheet.Range("C2").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Widen column C like the equivalent "result" column on the other sheets.
$newSheet.Columns.Item(3).ColumnWidth = 63.1666666666667

# Leave C2 selected / the new sheet active, as in the final workbook.
$newSheet.Range("C2").Select() | Out-Null
